$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 454789.9
$ws.Range("I9").Value = 268.9
$ws.Range("K9").Value = 268.9
$ws.Range("M9").Value = -99.89999999999998

$ws.Range("H32").Value = 1950
$ws.Range("J32").Value = 1950
$ws.Range("L32").Value = 1950
$ws.Range("N32").Value = -2602

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H125").Value = 1488.5
$ws.Range("I125").Value = 1488.5
$ws.Range("K125").Value = 13396.5
$ws.Range("M125").Value = -10936.5

$ws.Range("H132").Value = 456351
$ws.Range("I132").Value = 1891.5238
$ws.Range("K132").Value = 5674.5714
$ws.Range("M132").Value = -3144.5714

$ws.Range("H138").Value = 5166.2666
$ws.Range("I138").Value = 4049.25
$ws.Range("J138").Value = 5572.4546
$ws.Range("K138").Value = 12147.75
$ws.Range("L138").Value = 16717.3638
$ws.Range("M138").Value = -7007.75
$ws.Range("N138").Value = -26997.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18177
$ws.Range("I32").Value = 4376.067
$ws.Range("J32").Value = 59579.8
$ws.Range("K32").Value = 4376.067
$ws.Range("L32").Value = 59579.8
$ws.Range("M32").Value = -4089.067
$ws.Range("N32").Value = -60153.8

$ws.Range("H46").Value = 20898.5
$ws.Range("I46").Value = 14297.5
$ws.Range("K46").Value = 14297.5
$ws.Range("M46").Value = -13978.5

$ws.Range("H61").Value = 3468.5833
$ws.Range("J61").Value = 3889.5715
$ws.Range("L61").Value = 3889.5715
$ws.Range("N61").Value = -4313.5715

$ws.Range("H74").Value = 2282.6
$ws.Range("I74").Value = 2224.75
$ws.Range("J74").Value = 2514
$ws.Range("K74").Value = 2224.75
$ws.Range("L74").Value = 2514
$ws.Range("M74").Value = -1350.75
$ws.Range("N74").Value = -4262

$ws.Range("H77").Value = 2282.6
$ws.Range("I77").Value = 2224.75
$ws.Range("J77").Value = 2514
$ws.Range("K77").Value = 11123.75
$ws.Range("L77").Value = 12570
$ws.Range("M77").Value = -6755.75
$ws.Range("N77").Value = -21306

$ws.Range("H88").Value = 2214.3157
$ws.Range("J88").Value = 2392.3076
$ws.Range("L88").Value = 2392.3076
$ws.Range("N88").Value = -3204.3076

$ws.Range("H91").Value = 2214.3157
$ws.Range("J91").Value = 2392.3076
$ws.Range("L91").Value = 2392.3076
$ws.Range("N91").Value = -5200.3076

$ws.Range("H110").Value = 2715.7
$ws.Range("I110").Value = 2607.375
$ws.Range("K110").Value = 2607.375
$ws.Range("M110").Value = -562.375

$ws.Range("H122").Value = 2778.7942
$ws.Range("I122").Value = 2551.8262
$ws.Range("K122").Value = 7655.4786
$ws.Range("M122").Value = -5205.4786

$ws.Range("H136").Value = 3468.5833
$ws.Range("J136").Value = 3889.5715
$ws.Range("L136").Value = 11668.7145
$ws.Range("N136").Value = -16768.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2509.1667
$ws.Range("I86").Value = 2789
$ws.Range("K86").Value = 2789
$ws.Range("M86").Value = -1666

$ws.Range("H89").Value = 2509.1667
$ws.Range("I89").Value = 2789
$ws.Range("K89").Value = 13945
$ws.Range("M89").Value = -8329

$ws.Range("H99").Value = 2916.5
$ws.Range("I99").Value = 800
$ws.Range("J99").Value = 3974.75
$ws.Range("K99").Value = 800
$ws.Range("L99").Value = 3974.75
$ws.Range("M99").Value = 698
$ws.Range("N99").Value = -6970.75

$ws.Range("H105").Value = 5112
$ws.Range("J105").Value = 4999.8335
$ws.Range("L105").Value = 4999.8335
$ws.Range("N105").Value = -8493.833500000001

$ws.Range("H107").Value = 1830
$ws.Range("I107").Value = 1830
$ws.Range("K107").Value = 1830
$ws.Range("M107").Value = 90

$ws.Range("H134").Value = 2510.2856
$ws.Range("I134").Value = 1725.8422
$ws.Range("J134").Value = 4166.3335
$ws.Range("K134").Value = 5177.5266
$ws.Range("L134").Value = 12499.0005
$ws.Range("M134").Value = -2642.5266
$ws.Range("N134").Value = -17569.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 36203.5
$ws.Range("J93").Value = 62000
$ws.Range("L93").Value = 62000
$ws.Range("N93").Value = -65744

$ws.Range("H132").Value = 1783.7693
$ws.Range("I132").Value = 835.4545000000001
$ws.Range("K132").Value = 2506.3635
$ws.Range("M132").Value = 23.63649999999961

$ws.Range("H134").Value = 3079.5
$ws.Range("I134").Value = 2831.7576
$ws.Range("J134").Value = 5804.6665
$ws.Range("K134").Value = 8495.272799999999
$ws.Range("L134").Value = 17413.9995
$ws.Range("M134").Value = -5960.272799999999
$ws.Range("N134").Value = -22483.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1903.2
$ws.Range("I33").Value = 562.6
$ws.Range("K33").Value = 3375.6
$ws.Range("M33").Value = -3092.6

$ws.Range("H50").Value = 16664
$ws.Range("I50").Value = 19997
$ws.Range("J50").Value = 9998
$ws.Range("K50").Value = 59991
$ws.Range("L50").Value = 29994
$ws.Range("M50").Value = -59510
$ws.Range("N50").Value = -30956

$ws.Range("H53").Value = 16664
$ws.Range("I53").Value = 19997
$ws.Range("J53").Value = 9998
$ws.Range("K53").Value = 59991
$ws.Range("L53").Value = 29994
$ws.Range("M53").Value = -59510
$ws.Range("N53").Value = -30956

$ws.Range("H131").Value = 44425
$ws.Range("J131").Value = 1917.4736
$ws.Range("L131").Value = 5752.4208
$ws.Range("N131").Value = -15832.4208

$ws.Range("H137").Value = 2783
$ws.Range("I137").Value = 2783
$ws.Range("K137").Value = 8349
$ws.Range("M137").Value = -3249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11500000
$ws.Range("I11").Value = 6500000
$ws.Range("K11").Value = 6500000
$ws.Range("M11").Value = -6499861

$ws.Range("H92").Value = 14300.5
$ws.Range("J92").Value = 14300.5
$ws.Range("L92").Value = 14300.5
$ws.Range("N92").Value = -18044.5

$ws.Range("H102").Value = 38111.48
$ws.Range("I102").Value = 46207.1
$ws.Range("J102").Value = 5729
$ws.Range("K102").Value = 46207.1
$ws.Range("L102").Value = 5729
$ws.Range("M102").Value = -44585.1
$ws.Range("N102").Value = -8973

$ws.Range("H126").Value = 9999
$ws.Range("I126").Value = 9999
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 29997
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -27527
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 454.64285
$ws.Range("I16").Value = 454.64285
$ws.Range("K16").Value = 454.64285
$ws.Range("M16").Value = -284.64285

$ws.Range("H55").Value = 555.6667
$ws.Range("I55").Value = 683.6667
$ws.Range("K55").Value = 683.6667
$ws.Range("M55").Value = -510.6667

$ws.Range("H94").Value = 40165
$ws.Range("J94").Value = 40165
$ws.Range("L94").Value = 40165
$ws.Range("N94").Value = -41517

$ws.Range("H135").Value = 60382
$ws.Range("J135").Value = 60382
$ws.Range("L135").Value = 60382
$ws.Range("N135").Value = -70522

$ws.Range("H136").Value = 4804.1665
$ws.Range("J136").Value = 5589.636
$ws.Range("L136").Value = 16768.908
$ws.Range("N136").Value = -21868.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H39").Value = 29999
$ws.Range("J39").Value = 29999
$ws.Range("L39").Value = 29999
$ws.Range("N39").Value = -30825

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H132").Value = 1711.7576
$ws.Range("I132").Value = 1711.7576
$ws.Range("K132").Value = 5135.2728
$ws.Range("M132").Value = -2605.2728

$ws.Range("H137").Value = 64500
$ws.Range("J137").Value = 64500
$ws.Range("L137").Value = 64500
$ws.Range("N137").Value = -74700

Write-Host "Applied Raiden_Profits data updates"
